# Loan RBI, Variable Instalments
# Insert a new blank column before column N ("Late") on the
# "Repayment schedule" sheet, shifting Late/Date/Outstanding one column to
# the right, and make that sheet the active / selected tab.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("NewLoanInput")
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column) - pushes Late/Date/
# Outstanding from N:P to O:Q and leaves the new N column blank. The newly
# inserted column inherits the width of the column immediately to its
# left (M), matching Excel's native "Insert" behaviour.
$leftWidth = $wsSchedule.Columns("M").ColumnWidth
$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = $leftWidth

# Move the active selection / active tab from NewLoanInput to the
# Repayment schedule sheet.
$wsSchedule.Select()
$wsSchedule.Range("S8").Select()
